$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.178.28"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "1.854.57"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'238.26"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "'0.6915"
$ws.Range("E6").Value = "  -3.95%  "
$ws.Range("D7").Value = "'0.9993"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.07710"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").Value = "'0.3050"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "'23.26"
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("D11").Value = "'0.08059"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Value = "1.928.16"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "'0.7232"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").Value = "'5.211"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").Value = "'89.49"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").Value = "29.164.99"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").Value = "'5.748"
$ws.Range("E17").Value = "  -4.40%  "
$ws.Range("D18").Value = "'0.000007808"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "'234.66"
$ws.Range("E20").Value = "  -4.74%  "
$ws.Range("D21").Value = "'0.9988"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "2.101.90"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("D23").Value = "'0.9989"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "'7.428"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("D25").Value = "'161.64"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "'8.970"
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").Value = "'0.1436"
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("D28").Value = "'18.07"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").Value = "'1.962"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").Value = "'1.404"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").Value = "'4.515"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").Value = "'1.489"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").Value = "'4.024"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("D34").Value = "'0.05181"
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("D35").Value = "'1.185"
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("D36").Value = "'0.7052"
$ws.Range("E36").Value = "  -4.14%  "
$ws.Range("D37").Value = "'1.020"
$ws.Range("E37").Value = "  +2.03%  "
$ws.Range("D38").Value = "'2.671"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("D39").Value = "'0.01851"
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").Value = "'2.682"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "'0.9282"
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("D42").Value = "1.104.14"
$ws.Range("E42").Value = "  +6.32%  "
$ws.Range("D43").Value = "'5.937"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").Value = "'0.4282"
$ws.Range("E44").Value = "  -3.83%  "
$ws.Range("D45").Value = "'70.73"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").Value = "'0.9995"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "'101.80"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").Value = "'1.783"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "1.996.95"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").Value = "'9.175"
$ws.Range("E50").Value = "  -4.31%  "
$ws.Range("D51").Value = "'7.017"
$ws.Range("E51").Value = "  -5.83%  "
